$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting cell E8 with the new text
$ws.Range("E8").Value = "GIT UPDATE"

# Set the selection to reflect the last edited cell
$ws.Range("E8").Select()
